$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - rename to short codes
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize Spanish articles/prepositions (de/del/la/las/el/los/y) within names
$ws.Range("B7").Value = 'Pabellón De Arteaga'
$ws.Range("B8").Value = 'Rincón De Romos'
$ws.Range("B29").Value = 'Amatenango De La Frontera'
$ws.Range("B32").Value = 'Bejucal De Ocampo'
$ws.Range("B34").Value = 'Benemérito De Las Américas'
$ws.Range("B41").Value = 'Chiapa De Corzo'
$ws.Range("B47").Value = 'Comitán De Domínguez'
$ws.Range("B69").Value = 'Marqués De Comillas'
$ws.Range("B70").Value = 'Mazapa De Madero'
$ws.Range("B77").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B86").Value = 'Salto De Agua'
$ws.Range("B87").Value = 'San Cristóbal De Las Casas'
$ws.Range("B123").Value = 'Hidalgo Del Parral'
$ws.Range("B153").Value = 'San Juan De Sabinas'
$ws.Range("B165").Value = 'Villa De Álvarez'
$ws.Range("A167").Value = 'Ciudad De México'
$ws.Range("B171").Value = 'Cuajimalpa De Morelos'
$ws.Range("B194").Value = 'Nombre De Dios'
$ws.Range("B197").Value = 'Pánuco De Coronado'
$ws.Range("B203").Value = 'San Juan De Guadalupe'
$ws.Range("A212").Value = 'Estado De México'
$ws.Range("B212").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B215").Value = 'Almoloya De Alquisiras'
$ws.Range("B216").Value = 'Almoloya De Juárez'
$ws.Range("B222").Value = 'Atizapán De Zaragoza'
$ws.Range("B226").Value = 'Chapa De Mota'
$ws.Range("B229").Value = 'Coacalco De Berriozábal'
$ws.Range("B234").Value = 'Ecatepec De Morelos'
$ws.Range("B240").Value = 'Ixtapan De La Sal'
$ws.Range("B241").Value = 'Ixtapan Del Oro'
$ws.Range("B251").Value = 'Naucalpan De Juárez'
$ws.Range("B261").Value = 'San Antonio La Isla'
$ws.Range("B262").Value = 'San Felipe Del Progreso'
$ws.Range("B264").Value = 'San Simón De Guerero'
$ws.Range("B275").Value = 'Tenango Del Valle'
$ws.Range("B286").Value = 'Tlalnepantla De Baz'
$ws.Range("B291").Value = 'Valle De Bravo'
$ws.Range("B292").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B293").Value = 'Villa De Allende'
$ws.Range("B294").Value = 'Villa Del Carbón'
$ws.Range("A303").Value = 'Guanajuato'
$ws.Range("B306").Value = 'Apaseo El Alto'
$ws.Range("B307").Value = 'Apaseo El Grande'
$ws.Range("B315").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B319").Value = 'Jaral Del Progreso'
$ws.Range("B327").Value = 'Purísima Del Rincón'
$ws.Range("B331").Value = 'San Diego De La Unión'
$ws.Range("B333").Value = 'San Francisco Del Rincón'
$ws.Range("B335").Value = 'San Luis De La Paz'
$ws.Range("B336").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B338").Value = 'Silao De La Victoria'
$ws.Range("B343").Value = 'Valle De Santiago'
$ws.Range("B349").Value = 'Acapulco De Juárez'
$ws.Range("B352").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B353").Value = 'Alcozauca De Guerero'
$ws.Range("B357").Value = 'Atenango Del Río'
$ws.Range("B358").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B360").Value = 'Atoyac De Álvarez'
$ws.Range("B361").Value = 'Ayutla De Los Libres'
$ws.Range("B364").Value = 'Buenavista De Cuéllar'
$ws.Range("B365").Value = 'Chilapa De Álvarez'
$ws.Range("B366").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B367").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B372").Value = 'Coyuca De Benítez'
$ws.Range("B373").Value = 'Coyuca De Catalán'
$ws.Range("B376").Value = 'Cuetzala Del Progreso'
$ws.Range("B377").Value = 'Cutzamala De Pinzón'
$ws.Range("B383").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B384").Value = 'Iguala De La Independencia'
$ws.Range("B386").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B387").Value = 'Zihuatanejo De Azueta'
$ws.Range("B389").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B392").Value = 'Mártir De Cuilapan'
$ws.Range("B404").Value = 'Taxco De Alarcón'
$ws.Range("B406").Value = 'Técpan De Galeana'
$ws.Range("B408").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B410").Value = 'Tixtla De Guerero'
$ws.Range("B413").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B414").Value = 'Tlapa De Comonfort'
$ws.Range("B426").Value = 'Agua Blanca De Iturbide'
$ws.Range("B432").Value = 'Atotonilco De Tula'
$ws.Range("B433").Value = 'Atotonilco El Grande'
$ws.Range("B439").Value = 'Cuautepec De Hinojosa'
$ws.Range("B443").Value = 'Huasca De Ocampo'
$ws.Range("B447").Value = 'Huejutla De Reyes'
$ws.Range("B450").Value = 'Jacala De Ledezma'
$ws.Range("B455").Value = 'Mineral Del Chico'
$ws.Range("B456").Value = 'Mineral Del Monte'
$ws.Range("B457").Value = 'Mixquiahuala De Juárez'
$ws.Range("B458").Value = 'Molango De Escamilla'
$ws.Range("B460").Value = 'Nopala De Villagrán'
$ws.Range("B461").Value = 'Omitlán De Juárez'
$ws.Range("B462").Value = 'Pachuca De Soto'
$ws.Range("B465").Value = 'Progreso De Obregón'
$ws.Range("B470").Value = 'Santiago De Anaya'
$ws.Range("B471").Value = 'Santiago Tulantepec De Lugo Guerero'
$ws.Range("B475").Value = 'Tenango De Doria'
$ws.Range("B477").Value = 'Tepehuacán De Guerero'
$ws.Range("B478").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B481").Value = 'Tezontepec De Aldama'
$ws.Range("B487").Value = 'Tula De Allende'
$ws.Range("B488").Value = 'Tulancingo De Bravo'
$ws.Range("B489").Value = 'Villa De Tezontepec'
$ws.Range("B492").Value = 'Zacualtipán De Ángeles'
$ws.Range("B498").Value = 'Atotonilco El Alto'
$ws.Range("B499").Value = 'Autlán De Navarro'
$ws.Range("B513").Value = 'Encarnación De Díaz'
$ws.Range("B519").Value = 'Huejuquilla El Alto'
$ws.Range("B520").Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range("B521").Value = 'Ixtlahuacán Del Río'
$ws.Range("B528").Value = 'Lagos De Moreno'
$ws.Range("B532").Value = 'Ojuelos De Jalisco'
$ws.Range("B537").Value = 'San Diego De Alejandría'
$ws.Range("B538").Value = 'San Juan De Los Lagos'
$ws.Range("B539").Value = 'San Martín De Bolaños'
$ws.Range("B540").Value = 'San Miguel El Alto'
$ws.Range("B541").Value = 'San Sebastián Del Oeste'
$ws.Range("B544").Value = 'Tamazula De Gordiano'
$ws.Range("B546").Value = 'Techaluta De Montenegro'
$ws.Range("B549").Value = 'Teocuitatlán De Corona'
$ws.Range("B550").Value = 'Tepatitlán De Morelos'
$ws.Range("B552").Value = 'Tizapán El Alto'
$ws.Range("B553").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B559").Value = 'Unión De San Antonio'
$ws.Range("B560").Value = 'Unión De Tula'
$ws.Range("B563").Value = 'Yahualica De González Gallo'
$ws.Range("B564").Value = 'Zacoalco De Torres'
$ws.Range("B567").Value = 'Zapotitlán De Vadillo'
$ws.Range("B568").Value = 'Zapotlán El Grande'
$ws.Range("B591").Value = 'Cojumatlán De Régules'
$ws.Range("B649").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B672").Value = 'Coatlán Del Río'
$ws.Range("B680").Value = 'Jonacatepec De Leandro Valle'
$ws.Range("B684").Value = 'Puente De Ixtla'
$ws.Range("B689").Value = 'Tetela Del Volcán'
$ws.Range("B690").Value = 'Tlaltizapán De Zapata'
$ws.Range("B697").Value = 'Zacualpan De Amilpas'
$ws.Range("B700").Value = 'Amatlán De Cañas'
$ws.Range("B701").Value = 'Bahía De Banderas'
$ws.Range("B703").Value = 'Ixtlán Del Río'
$ws.Range("B710").Value = 'Santa María Del Oro'
$ws.Range("B734").Value = 'Mier Y Noriega'
$ws.Range("B735").Value = 'Montemorelos'
$ws.Range("B739").Value = 'San Nicolás De Los Garza'
$ws.Range("B743").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B748").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B751").Value = 'Fresnillo De Trujano'
$ws.Range("B752").Value = 'Guevea De Humboldt'
$ws.Range("B753").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B754").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B755").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B756").Value = 'Huautla De Jiménez'
$ws.Range("B757").Value = 'Ixtlán De Juárez'
$ws.Range("B758").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B766").Value = 'Mariscala De Juárez'
$ws.Range("B767").Value = 'Mártires De Tacubaya'
$ws.Range("B769").Value = 'Mazatlán Villa De Flores'
$ws.Range("B770").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B771").Value = 'Mixistlán De La Reforma'
$ws.Range("B773").Value = 'Nejapa De Madero'
$ws.Range("B774").Value = 'Oaxaca De Juárez'
$ws.Range("B775").Value = 'Ocotlán De Morelos'
$ws.Range("B776").Value = 'Pinotepa De Don Luis'
$ws.Range("B778").Value = 'Putla Villa De Guerero'
$ws.Range("B779").Value = 'Reforma De Pineda'
$ws.Range("B800").Value = 'San Felipe Jalapa De Díaz'
$ws.Range("B820").Value = 'San Juan Bautista Lo De Soto'
$ws.Range("B850").Value = 'San Martín De Los Cansecos'
$ws.Range("B863").Value = 'San Miguel Del Puerto'
$ws.Range("B864").Value = 'San Miguel El Grande'
$ws.Range("B877").Value = 'San Pablo Villa De Mitla'
$ws.Range("B883").Value = 'San Pedro El Alto'
$ws.Range("B900").Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range("B901").Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Range("B917").Value = 'Santa Cruz Tacache De Mina'
$ws.Range("B921").Value = 'Santa Inés Del Monte'
$ws.Range("B931").Value = 'Santa María Del Rosario'
$ws.Range("B932").Value = 'Santa María Del Tule'
$ws.Range("B939").Value = 'Santa María Jalapa Del Marqués'
$ws.Range("B963").Value = 'Santiago Del Río'
$ws.Range("B990").Value = 'Santo Domingo De Morelos'
$ws.Range("B1002").Value = 'Tanetze De Zaragoza'
$ws.Range("B1003").Value = 'Tataltepec De Valdés'
$ws.Range("B1004").Value = 'Teotitlán De Flores Magón'
$ws.Range("B1005").Value = 'Tezoatlán De Segura Y Luna'
$ws.Range("B1006").Value = 'Tlacolula De Matamoros'
$ws.Range("B1008").Value = 'Villa De Etla'
$ws.Range("B1009").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B1011").Value = 'Villa Sola De Vega'
$ws.Range("B1012").Value = 'Zapotitlán Del Río'
$ws.Range("B1034").Value = 'Ayotoxco De Guerero'
$ws.Range("B1036").Value = 'Chalchicomula De Sesma'
$ws.Range("B1044").Value = 'Chila De La Sal'
$ws.Range("B1054").Value = 'Cuayuca De Andrade'
$ws.Range("B1066").Value = 'Huehuetlán El Chico'
$ws.Range("B1069").Value = 'Ixcamilpa De Guerero'
$ws.Range("B1071").Value = 'Izúcar De Matamoros'
$ws.Range("B1084").Value = 'Palmar De Bravo'
$ws.Range("B1105").Value = 'San Nicolás De Los Ranchos'
$ws.Range("B1107").Value = 'San Salvador El Seco'
$ws.Range("B1108").Value = 'San Salvador El Verde'
$ws.Range("B1112").Value = 'Tecali De Herrera'
$ws.Range("B1118").Value = 'Tepanco De López'
$ws.Range("B1119").Value = 'Tepango De Rodríguez'
$ws.Range("B1120").Value = 'Tepatlaxco De Hidalgo'
$ws.Range("B1123").Value = 'Tepexi De Rodríguez'
$ws.Range("B1124").Value = 'Tepeyahualco De Cuauhtémoc'
$ws.Range("B1125").Value = 'Tetela De Ocampo'
$ws.Range("B1130").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B1145").Value = 'Xayacatlán De Bravo'
$ws.Range("B1160").Value = 'Amealco De Bonfil'
$ws.Range("B1162").Value = 'Cadereyta De Montes'
$ws.Range("B1168").Value = 'Jalpan De Serra'
$ws.Range("B1169").Value = 'Landa De Matamoros'
$ws.Range("B1172").Value = 'Pinal De Amoles'
$ws.Range("B1175").Value = 'San Juan Del Río'
$ws.Range("B1186").Value = 'Axtla De Terrazas'
$ws.Range("B1192").Value = 'Ciudad Del Maíz'
$ws.Range("B1201").Value = 'Mexquitic De Carmona'
$ws.Range("B1207").Value = 'San Ciro De Acosta'
$ws.Range("B1212").Value = 'Santa María Del Río'
$ws.Range("B1214").Value = 'Soledad De Graciano Sánchez'
$ws.Range("B1224").Value = 'Villa De Arista'
$ws.Range("B1225").Value = 'Villa De Arriaga'
$ws.Range("B1226").Value = 'Villa De Guadalupe'
$ws.Range("B1227").Value = 'Villa De Ramos'
$ws.Range("B1228").Value = 'Villa De Reyes'
$ws.Range("B1262").Value = 'Nacozari De García'
$ws.Range("B1275").Value = 'Jalpa De Méndez'
$ws.Range("B1311").Value = 'Soto La Marina'
$ws.Range("B1324").Value = 'Contla De Juan Cuamatzi'
$ws.Range("B1327").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B1330").Value = 'San Pablo Del Monte'
$ws.Range("B1334").Value = 'Tepetitla De Lardizábal'
$ws.Range("B1337").Value = 'Tetla De La Solidaridad'
$ws.Range("B1352").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B1356").Value = 'Amatlán De Los Reyes'
$ws.Range("B1366").Value = 'Boca Del Río'
$ws.Range("B1371").Value = 'Castillo De Teayo'
$ws.Range("B1373").Value = 'Cazones De Herrera'
$ws.Range("B1389").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1405").Value = 'Hueyapan De Ocampo'
$ws.Range("B1406").Value = 'Ignacio De La Llave'
$ws.Range("B1410").Value = 'Ixhuatlán De Madero'
$ws.Range("B1411").Value = 'Ixhuatlán Del Café'
$ws.Range("B1412").Value = 'Ixhuatlán Del Sureste'
$ws.Range("B1421").Value = 'Juchique De Ferrer'
$ws.Range("B1424").Value = 'Las Vigas De Ramírez'
$ws.Range("B1425").Value = 'Lerdo De Tejada'
$ws.Range("B1428").Value = 'Martínez De La Torre'
$ws.Range("B1431").Value = 'Medellín De Bravo'
$ws.Range("B1435").Value = 'Mixtla De Altamirano'
$ws.Range("B1437").Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range("B1446").Value = 'Ozuluama De Mascareñas'
$ws.Range("B1450").Value = 'Paso De Ovejas'
$ws.Range("B1451").Value = 'Paso Del Macho'
$ws.Range("B1455").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1463").Value = 'Sayula De Alemán'
$ws.Range("B1466").Value = 'Soledad De Doblado'
$ws.Range("B1471").Value = 'Tatahuicapan De Juárez'
$ws.Range("B1501").Value = 'Vega De Alatorre'
$ws.Range("B1511").Value = 'Zozocolco De Hidalgo'
$ws.Range("B1525").Value = 'Cañitas De Felipe Pescador'
$ws.Range("B1527").Value = 'Concepción Del Oro'
$ws.Range("B1536").Value = 'Jiménez Del Teul'
$ws.Range("B1543").Value = 'Mezquital Del Oro'
$ws.Range("B1547").Value = 'Nochistlán De Mejía'
$ws.Range("B1548").Value = 'Noria De Ángeles'
$ws.Range("B1555").Value = 'Teúl De González Ortega'
$ws.Range("B1556").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1559").Value = 'Villa De Cos'

# Floating point precision corrections
$ws.Range("D6").Value = 0.0000921871398939848
$ws.Range("D30").Value = 0.0000921871398939848
$ws.Range("D35").Value = 0.0000921871398939848
$ws.Range("D40").Value = 0.0000921871398939848
$ws.Range("D41").Value = 0.0000921871398939848
$ws.Range("D57").Value = 0.0000921871398939848
$ws.Range("D59").Value = 0.0000921871398939848
$ws.Range("D60").Value = 0.0000921871398939848
$ws.Range("D61").Value = 0.0000921871398939848
$ws.Range("D76").Value = 0.0000921871398939848
$ws.Range("D80").Value = 0.0000921871398939848
$ws.Range("D95").Value = 0.0000921871398939848
$ws.Range("D98").Value = 0.0000921871398939848
$ws.Range("D104").Value = 0.0000921871398939848
$ws.Range("D111").Value = 0.0000921871398939848
$ws.Range("D115").Value = 0.0000921871398939848
$ws.Range("D122").Value = 0.0000921871398939848
$ws.Range("D133").Value = 0.0000921871398939848
$ws.Range("D135").Value = 0.0000921871398939848
$ws.Range("D139").Value = 0.0000921871398939848
$ws.Range("D157").Value = 0.0000921871398939848
$ws.Range("D179").Value = 0.0009679649688868404
$ws.Range("D194").Value = 0.0000921871398939848
$ws.Range("D196").Value = 0.0000921871398939848
$ws.Range("D197").Value = 0.0000921871398939848
$ws.Range("D198").Value = 0.0000921871398939848
$ws.Range("D207").Value = 0.0000921871398939848
$ws.Range("D220").Value = 0.0000921871398939848
$ws.Range("D229").Value = 0.0000921871398939848
$ws.Range("D233").Value = 0.0000921871398939848
$ws.Range("D237").Value = 0.0000921871398939848
$ws.Range("D238").Value = 0.0000921871398939848
$ws.Range("D255").Value = 0.0000921871398939848
$ws.Range("D256").Value = 0.0009679649688868404
$ws.Range("D263").Value = 0.0000921871398939848
$ws.Range("D270").Value = 0.0000921871398939848
$ws.Range("D277").Value = 0.0000921871398939848
$ws.Range("D279").Value = 0.0000921871398939848
$ws.Range("D280").Value = 0.0000921871398939848
$ws.Range("D283").Value = 0.0000921871398939848
$ws.Range("D294").Value = 0.0000921871398939848
$ws.Range("D309").Value = 0.009725743258815396
$ws.Range("D323").Value = 0.0000921871398939848
$ws.Range("D350").Value = 0.0000921871398939848
$ws.Range("D380").Value = 0.0000921871398939848
$ws.Range("D415").Value = 0.0009679649688868404
$ws.Range("D426").Value = 0.0009679649688868404
$ws.Range("D427").Value = 0.0000921871398939848
$ws.Range("D428").Value = 0.0009679649688868404
$ws.Range("D431").Value = 0.0000921871398939848
$ws.Range("D440").Value = 0.0000921871398939848
$ws.Range("D451").Value = 0.0000921871398939848
$ws.Range("D468").Value = 0.0000921871398939848
$ws.Range("D476").Value = 0.0000921871398939848
$ws.Range("D490").Value = 0.0000921871398939848
$ws.Range("D506").Value = 0.0000921871398939848
$ws.Range("D511").Value = 0.0000921871398939848
$ws.Range("D512").Value = 0.0000921871398939848
$ws.Range("D523").Value = 0.0000921871398939848
$ws.Range("D526").Value = 0.0000921871398939848
$ws.Range("D532").Value = 0.0000921871398939848
$ws.Range("D537").Value = 0.0000921871398939848
$ws.Range("D538").Value = 0.0000921871398939848
$ws.Range("D540").Value = 0.0000921871398939848
$ws.Range("D542").Value = 0.0000921871398939848
$ws.Range("D546").Value = 0.0000921871398939848
$ws.Range("D551").Value = 0.0000921871398939848
$ws.Range("D560").Value = 0.0000921871398939848
$ws.Range("D566").Value = 0.0000921871398939848
$ws.Range("D568").Value = 0.0000921871398939848
$ws.Range("D575").Value = 0.0009679649688868404
$ws.Range("D579").Value = 0.0000921871398939848
$ws.Range("D595").Value = 0.0000921871398939848
$ws.Range("D606").Value = 0.0000921871398939848
$ws.Range("D612").Value = 0.0000921871398939848
$ws.Range("D627").Value = 0.0000921871398939848
$ws.Range("D631").Value = 0.0000921871398939848
$ws.Range("D647").Value = 0.0000921871398939848
$ws.Range("D648").Value = 0.0000921871398939848
$ws.Range("D653").Value = 0.0009679649688868404
$ws.Range("D655").Value = 0.0009679649688868404
$ws.Range("D676").Value = 0.0000921871398939848
$ws.Range("D681").Value = 0.0000921871398939848
$ws.Range("D700").Value = 0.0000921871398939848
$ws.Range("D715").Value = 0.0000921871398939848
$ws.Range("D717").Value = 0.0000921871398939848
$ws.Range("D721").Value = 0.0000921871398939848
$ws.Range("D724").Value = 0.0000921871398939848
$ws.Range("D732").Value = 0.0000921871398939848
$ws.Range("D738").Value = 0.0000921871398939848
$ws.Range("D751").Value = 0.0000921871398939848
$ws.Range("D762").Value = 0.0000921871398939848
$ws.Range("D764").Value = 0.0000921871398939848
$ws.Range("D776").Value = 0.0000921871398939848
$ws.Range("D777").Value = 0.0000921871398939848
$ws.Range("D786").Value = 0.0000921871398939848
$ws.Range("D788").Value = 0.0000921871398939848
$ws.Range("D789").Value = 0.0000921871398939848
$ws.Range("D791").Value = 0.0000921871398939848
$ws.Range("D793").Value = 0.0000921871398939848
$ws.Range("D802").Value = 0.0000921871398939848
$ws.Range("D812").Value = 0.0000921871398939848
$ws.Range("D821").Value = 0.0000921871398939848
$ws.Range("D825").Value = 0.0000921871398939848
$ws.Range("D830").Value = 0.0000921871398939848
$ws.Range("D840").Value = 0.0000921871398939848
$ws.Range("D860").Value = 0.0000921871398939848
$ws.Range("D863").Value = 0.0000921871398939848
$ws.Range("D864").Value = 0.0000921871398939848
$ws.Range("D867").Value = 0.0000921871398939848
$ws.Range("D868").Value = 0.0000921871398939848
$ws.Range("D871").Value = 0.0000921871398939848
$ws.Range("D876").Value = 0.0000921871398939848
$ws.Range("D881").Value = 0.0000921871398939848
$ws.Range("D892").Value = 0.0000921871398939848
$ws.Range("D896").Value = 0.0000921871398939848
$ws.Range("D898").Value = 0.0000921871398939848
$ws.Range("D915").Value = 0.0000921871398939848
$ws.Range("D920").Value = 0.0000921871398939848
$ws.Range("D927").Value = 0.0000921871398939848
$ws.Range("D930").Value = 0.0000921871398939848
$ws.Range("D933").Value = 0.0000921871398939848
$ws.Range("D943").Value = 0.0000921871398939848
$ws.Range("D944").Value = 0.0000921871398939848
$ws.Range("D948").Value = 0.0000921871398939848
$ws.Range("D952").Value = 0.0000921871398939848
$ws.Range("D964").Value = 0.0000921871398939848
$ws.Range("D974").Value = 0.0000921871398939848
$ws.Range("D975").Value = 0.0000921871398939848
$ws.Range("D983").Value = 0.0000921871398939848
$ws.Range("D993").Value = 0.0000921871398939848
$ws.Range("D1001").Value = 0.0000921871398939848
$ws.Range("D1008").Value = 0.0000921871398939848
$ws.Range("D1009").Value = 0.0009679649688868404
$ws.Range("D1013").Value = 0.0000921871398939848
$ws.Range("D1020").Value = 0.0000921871398939848
$ws.Range("D1021").Value = 0.0000921871398939848
$ws.Range("D1022").Value = 0.0000921871398939848
$ws.Range("D1025").Value = 0.0000921871398939848
$ws.Range("D1030").Value = 0.0009679649688868404
$ws.Range("D1033").Value = 0.0000921871398939848
$ws.Range("D1036").Value = 0.0000921871398939848
$ws.Range("D1038").Value = 0.0000921871398939848
$ws.Range("D1047").Value = 0.0000921871398939848
$ws.Range("D1053").Value = 0.0000921871398939848
$ws.Range("D1054").Value = 0.0000921871398939848
$ws.Range("D1064").Value = 0.0000921871398939848
$ws.Range("D1074").Value = 0.0000921871398939848
$ws.Range("D1079").Value = 0.0000921871398939848
$ws.Range("D1098").Value = 0.0000921871398939848
$ws.Range("D1104").Value = 0.0000921871398939848
$ws.Range("D1107").Value = 0.0000921871398939848
$ws.Range("D1110").Value = 0.0000921871398939848
$ws.Range("D1112").Value = 0.0000921871398939848
$ws.Range("D1113").Value = 0.0009679649688868404
$ws.Range("D1118").Value = 0.0000921871398939848
$ws.Range("D1120").Value = 0.0000921871398939848
$ws.Range("D1121").Value = 0.0000921871398939848
$ws.Range("D1123").Value = 0.0000921871398939848
$ws.Range("D1125").Value = 0.0000921871398939848
$ws.Range("D1127").Value = 0.0000921871398939848
$ws.Range("D1141").Value = 0.0000921871398939848
$ws.Range("D1147").Value = 0.0000921871398939848
$ws.Range("D1153").Value = 0.0000921871398939848
$ws.Range("D1157").Value = 0.0000921871398939848
$ws.Range("D1178").Value = 0.0000921871398939848
$ws.Range("D1189").Value = 0.0000921871398939848
$ws.Range("D1191").Value = 0.0000921871398939848
$ws.Range("D1218").Value = 0.0000921871398939848
$ws.Range("D1225").Value = 0.0000921871398939848
$ws.Range("D1230").Value = 0.0000921871398939848
$ws.Range("D1236").Value = 0.0000921871398939848
$ws.Range("D1239").Value = 0.0009679649688868404
$ws.Range("D1242").Value = 0.0000921871398939848
$ws.Range("D1245").Value = 0.0000921871398939848
$ws.Range("D1248").Value = 0.0000921871398939848
$ws.Range("D1253").Value = 0.0000921871398939848
$ws.Range("D1261").Value = 0.0000921871398939848
$ws.Range("D1272").Value = 0.0000921871398939848
$ws.Range("D1276").Value = 0.0000921871398939848
$ws.Range("D1278").Value = 0.0000921871398939848
$ws.Range("D1288").Value = 0.0000921871398939848
$ws.Range("D1297").Value = 0.0000921871398939848
$ws.Range("D1299").Value = 0.0000921871398939848
$ws.Range("D1306").Value = 0.0000921871398939848
$ws.Range("D1307").Value = 0.0009679649688868404
$ws.Range("D1316").Value = 0.0000921871398939848
$ws.Range("D1319").Value = 0.0000921871398939848
$ws.Range("D1321").Value = 0.0000921871398939848
$ws.Range("D1324").Value = 0.0000921871398939848
$ws.Range("D1325").Value = 0.0000921871398939848
$ws.Range("D1337").Value = 0.0000921871398939848
$ws.Range("D1342").Value = 0.0000921871398939848
$ws.Range("D1346").Value = 0.0000921871398939848
$ws.Range("D1350").Value = 0.0000921871398939848
$ws.Range("D1360").Value = 0.0000921871398939848
$ws.Range("D1368").Value = 0.0000921871398939848
$ws.Range("D1371").Value = 0.0000921871398939848
$ws.Range("D1379").Value = 0.0000921871398939848
$ws.Range("D1380").Value = 0.0000921871398939848
$ws.Range("D1386").Value = 0.0000921871398939848
$ws.Range("D1387").Value = 0.0000921871398939848
$ws.Range("D1395").Value = 0.0000921871398939848
$ws.Range("D1409").Value = 0.0000921871398939848
$ws.Range("D1418").Value = 0.0000921871398939848
$ws.Range("D1431").Value = 0.0000921871398939848
$ws.Range("D1432").Value = 0.0000921871398939848
$ws.Range("D1453").Value = 0.0000921871398939848
$ws.Range("D1455").Value = 0.0009679649688868404
$ws.Range("D1459").Value = 0.0000921871398939848
$ws.Range("D1482").Value = 0.0009679649688868404
$ws.Range("D1484").Value = 0.0009679649688868404
$ws.Range("D1486").Value = 0.0000921871398939848
$ws.Range("D1488").Value = 0.0000921871398939848
$ws.Range("D1498").Value = 0.0000921871398939848
$ws.Range("D1518").Value = 0.0000921871398939848
$ws.Range("D1520").Value = 0.0000921871398939848
$ws.Range("D1526").Value = 0.0000921871398939848
$ws.Range("D1535").Value = 0.0000921871398939848
$ws.Range("D1541").Value = 0.0000921871398939848
$ws.Range("D1544").Value = 0.0000921871398939848
$ws.Range("D1545").Value = 0.0000921871398939848
$ws.Range("D1546").Value = 0.0000921871398939848
$ws.Range("D1548").Value = 0.0000921871398939848
$ws.Range("D1562").Value = 0.0000921871398939848

# Delete footer/metadata rows 1568-1572
$ws.Range("A1568:D1572").EntireRow.Delete()